$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each changed cell: temporarily force Text number format so
# numeric-looking strings (prices formatted with "." thousands
# separators, trailing zeros, percentages, etc.) are written
# verbatim as text instead of being reinterpreted as numbers,
# then restore the default "Normal" style so the cell style stays
# unchanged from the source workbook.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '29.288.71'
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  +0.54%  '
$ws.Range('E2').Style = "Normal"

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.934.60'
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.49%  '
$ws.Range('E3').Style = "Normal"

$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E4').Style = "Normal"

$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -0.03%  '
$ws.Range('E5').Style = "Normal"

$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.9998'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('E6').Style = "Normal"

$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('E7').Style = "Normal"

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3875'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -0.26%  '
$ws.Range('E8').Style = "Normal"

$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('E9').Style = "Normal"

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07821'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -0.46%  '
$ws.Range('E10').Style = "Normal"

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.9745'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -1.59%  '
$ws.Range('E11').Style = "Normal"

$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +3.21%  '
$ws.Range('E12').Style = "Normal"

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '1.940.38'
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +3.33%  '
$ws.Range('E13').Style = "Normal"

$ws.Range('B14').NumberFormat = "@"
$ws.Range('B14').Value = 'Chainlink'
$ws.Range('B14').Style = "Normal"
$ws.Range('C14').NumberFormat = "@"
$ws.Range('C14').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('C14').Style = "Normal"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '7.083'
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +0.60%  '
$ws.Range('E14').Style = "Normal"

$ws.Range('B15').NumberFormat = "@"
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('B15').Style = "Normal"
$ws.Range('C15').NumberFormat = "@"
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('C15').Style = "Normal"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '5.785'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +0.63%  '
$ws.Range('E15').Style = "Normal"

$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.07083'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +0.87%  '
$ws.Range('E16').Style = "Normal"

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '86.76'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -1.44%  '
$ws.Range('E17').Style = "Normal"

$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -0.02%  '
$ws.Range('E18').Style = "Normal"

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.000009737'
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  -1.90%  '
$ws.Range('E19').Style = "Normal"

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '16.98'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -0.48%  '
$ws.Range('E20').Style = "Normal"

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  +0.08%  '
$ws.Range('E21').Style = "Normal"

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '29.268.19'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.46%  '
$ws.Range('E22').Style = "Normal"

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '5.467'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +2.72%  '
$ws.Range('E23').Style = "Normal"

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '11.05'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -0.54%  '
$ws.Range('E24').Style = "Normal"

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.171.99'
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +2.87%  '
$ws.Range('E25').Style = "Normal"

$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +0.31%  '
$ws.Range('E26').Style = "Normal"

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '157.57'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +1.00%  '
$ws.Range('E27').Style = "Normal"

$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -0.19%  '
$ws.Range('E28').Style = "Normal"

$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.775'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -2.30%  '
$ws.Range('E29').Style = "Normal"

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '118.90'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +0.18%  '
$ws.Range('E30').Style = "Normal"

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.835'
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -2.16%  '
$ws.Range('E31').Style = "Normal"

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '0.09344'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +0.31%  '
$ws.Range('E32').Style = "Normal"

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.8642'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -3.47%  '
$ws.Range('E33').Style = "Normal"

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '5.164'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -1.19%  '
$ws.Range('E34').Style = "Normal"

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.300'
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -1.50%  '
$ws.Range('E35').Style = "Normal"

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '3.083'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -1.86%  '
$ws.Range('E36').Style = "Normal"

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.05779'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -0.01%  '
$ws.Range('E37').Style = "Normal"

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.156'
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -1.09%  '
$ws.Range('E38').Style = "Normal"

$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.02079'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('E39').Style = "Normal"

$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '7.638'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -0.15%  '
$ws.Range('E40').Style = "Normal"

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.5656'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -0.78%  '
$ws.Range('E41').Style = "Normal"

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.000003103'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +11.85%  '
$ws.Range('E42').Style = "Normal"

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.1779'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  -1.62%  '
$ws.Range('E43').Style = "Normal"

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '9.394'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -3.39%  '
$ws.Range('E44').Style = "Normal"

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '2.710'
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  +6.39%  '
$ws.Range('E45').Style = "Normal"

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.5257'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -1.81%  '
$ws.Range('E46').Style = "Normal"

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '11.41'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -3.75%  '
$ws.Range('E47').Style = "Normal"

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.06857'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -1.78%  '
$ws.Range('E48').Style = "Normal"

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '2.067'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -4.47%  '
$ws.Range('E49').Style = "Normal"

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.809'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -1.61%  '
$ws.Range('E50').Style = "Normal"

$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '111.26'
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -1.55%  '
$ws.Range('E51').Style = "Normal"
